{"js": "// Insert the missing sentence about notifying the Taxi Driver into the\n// \"Exceptions\" bullet that talks about redistributing taxis between\n// queues. The sentence is added right after \"...the taxis in the city.\"\n// and before the trailing \" (Chiedere a Sofia)\" note.\n\nconst body = context.document.body;\n\nconst searchResults = body.search(\"the taxis in the city.\", { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find anchor text \"the taxis in the city.\" in the document.');\n}\n\nconst anchor = searchResults.items[0];\nanchor.insertText(\n  \" Then send a notification to the Taxi Driver and he moves to the new zone.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Insert the missing sentence about notifying the Taxi Driver into the\n# \"Exceptions\" bullet that talks about redistributing taxis between\n# queues. The sentence is added right after \"...the taxis in the city.\"\n# and before the trailing \" (Chiedere a Sofia)\" note.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"the taxis in the city.\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw 'Could not find anchor text \"the taxis in the city.\" in the document.'\n}\n\n# Collapse the found range to its end point so the insertion lands right\n# after the period, then add the new sentence there.\n$rng.Collapse([Microsoft.Office.Interop.Word.WdCollapseDirection]::wdCollapseEnd)\n$rng.InsertAfter(\" Then send a notification to the Taxi Driver and he moves to the new zone.\")\n"}
